$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A1:U60")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
try {
$tbl.TableStyle = 0
"set int ok"
} catch {
"err: " + $_.Exception.Message
}
